$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy formatting (bold, border,
# centered) from the existing header cell H1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-27.
$data = @(
    @(3, 6),
    @(4, 6),
    @(4, 5),
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 2),
    @(2, 2)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
